$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Text tweak: "Describe what plans ... and splice ..." ->
#    "Describe what the plans ... next increment are and splice ..."
# ---------------------------------------------------------------------------
$oldSentence = "Describe what plans for our projects next increment and splice all team members videos together."
$newSentence = "Describe what the plans for our projects next increment are and splice all team members videos together."

$d.Content.Find.Execute($oldSentence, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $newSentence, 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark: it used to live alone in its own empty
#    paragraph right after the "Link to video" line; it now belongs right
#    after "...next increment are" (and before " and splice...") inside the
#    sentence we just edited. Adding a bookmark with the same name moves it
#    (Word bookmark names are unique), and the now-empty paragraph that used
#    to hold it gets deleted.
# ---------------------------------------------------------------------------

# Locate the paragraph that now contains the edited sentence.
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Contains("next increment are and splice")) {
        $targetPara = $p
        break
    }
}

if ($targetPara -ne $null) {
    $searchRange = $targetPara.Range.Duplicate
    $searchRange.Find.Execute("next increment are", $true, $false, $false, $false, $false, `
                               $true, 1, $false, "", 0) | Out-Null

    $bmRange = $d.Range($searchRange.End, $searchRange.End)
    $d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
}

# Remove the paragraph that used to contain the bookmark by itself: it sits
# directly after the "Link to video" paragraph and, once the bookmark has
# moved out of it, is an empty paragraph.
$linkParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.Contains("Link to video")) {
        $linkParaIndex = $i
        break
    }
}

if ($linkParaIndex -gt 0) {
    $afterLinkPara = $d.Paragraphs.Item($linkParaIndex + 1)
    if ($afterLinkPara.Range.Text.Trim() -eq "") {
        $afterLinkPara.Range.Delete() | Out-Null
    }
}
